# Update "想去人数" (column F) values per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1370
$ws.Cells.Item(5, 6).Value = 5638
$ws.Cells.Item(6, 6).Value = 416
$ws.Cells.Item(7, 6).Value = 1037
$ws.Cells.Item(8, 6).Value = 2753
$ws.Cells.Item(9, 6).Value = 6404
$ws.Cells.Item(10, 6).Value = 175
$ws.Cells.Item(11, 6).Value = 1233
$ws.Cells.Item(12, 6).Value = 717
$ws.Cells.Item(13, 6).Value = 84
$ws.Cells.Item(14, 6).Value = 6
$ws.Cells.Item(15, 6).Value = 1100
$ws.Cells.Item(17, 6).Value = 74
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(19, 6).Value = 144
$ws.Cells.Item(21, 6).Value = 894
$ws.Cells.Item(25, 6).Value = 1128
$ws.Cells.Item(28, 6).Value = 227
$ws.Cells.Item(30, 6).Value = 224
$ws.Cells.Item(32, 6).Value = 44

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(23, 6).Value = 87
$ws.Cells.Item(27, 6).Value = 606
$ws.Cells.Item(35, 6).Value = 121
$ws.Cells.Item(37, 6).Value = 47

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 694
$ws.Cells.Item(5, 6).Value = 816
$ws.Cells.Item(7, 6).Value = 277

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 694
$ws.Cells.Item(4, 6).Value = 1370
$ws.Cells.Item(6, 6).Value = 816
$ws.Cells.Item(10, 6).Value = 277
$ws.Cells.Item(11, 6).Value = 277
$ws.Cells.Item(14, 6).Value = 5638
$ws.Cells.Item(15, 6).Value = 416
$ws.Cells.Item(16, 6).Value = 1037
$ws.Cells.Item(17, 6).Value = 2753
$ws.Cells.Item(19, 6).Value = 6404
$ws.Cells.Item(21, 6).Value = 175
$ws.Cells.Item(22, 6).Value = 1233
$ws.Cells.Item(25, 6).Value = 717
$ws.Cells.Item(26, 6).Value = 84
$ws.Cells.Item(27, 6).Value = 1100
$ws.Cells.Item(29, 6).Value = 74
$ws.Cells.Item(30, 6).Value = 144
$ws.Cells.Item(32, 6).Value = 894
$ws.Cells.Item(33, 6).Value = 87
$ws.Cells.Item(35, 6).Value = 1128
$ws.Cells.Item(38, 6).Value = 17
$ws.Cells.Item(41, 6).Value = 227
$ws.Cells.Item(44, 6).Value = 224
$ws.Cells.Item(47, 6).Value = 121
$ws.Cells.Item(50, 6).Value = 47
